# fix: correção na 2.3E
# Re-solves the "2.3A - 4" and "2.3E - 8" sheets with corrected inputs and
# refreshes the Solver bookkeeping (defined names + final selections).

$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("2.3A - 4")
$wsE = $wb.Worksheets.Item("2.3E - 8")
$wsC = $wb.Worksheets.Item("2.3C - 3")
$ws9A = $wb.Worksheets.Item("9.1A - 7")
$ws9B = $wb.Worksheets.Item("9.1B - 3")

# ---------------------------------------------------------------------
# "2.3A - 4": updated optimal solution (B19:F19) and the K9 constraint
# formula now also accounts for F19 (0.034*F19 term added).
# ---------------------------------------------------------------------
$wsA.Range("B19").Value = 131
$wsA.Range("C19").Value = 129
$wsA.Range("D19").Value = 207
$wsA.Range("E19").Value = 301
$wsA.Range("F19").Value = 34

$wsA.Range("K9").Formula = "=0.023*B19+0.034*C19+0.046*D19+0.023*E19 + 0.034*F19"

# Solver options on this sheet were tightened.
$nPre = $wsA.Names.Item("2.3A - 4!solver_pre")
$nPre.RefersTo = "=0.0000001"
$nTol = $wsA.Names.Item("2.3A - 4!solver_tol")
$nTol.RefersTo = "=0"
$nVal = $wsA.Names.Item("2.3A - 4!solver_val")
$nVal.RefersTo = "=8290"

# ---------------------------------------------------------------------
# "2.3E - 8": corrected coefficients (B8/C8 swapped) and the re-solved
# production quantities (B11:G11).
# ---------------------------------------------------------------------
$wsE.Range("B8").Value = 0
$wsE.Range("C8").Value = -12

$wsE.Range("B11").Value = 275.00000000000205
$wsE.Range("C11").Value = 874.99999999999784
$wsE.Range("D11").Value = 61.111111111111619
$wsE.Range("E11").Value = 388.8888888888884
$wsE.Range("F11").Value = 213.88888888889039
$wsE.Range("G11").Value = 486.11111111110961

$excel.Calculate()

# Re-assert these comparison formulas verbatim so the engine re-evaluates
# them against the freshly recalculated H column (they sit on cells that
# are fed by the MMULT array-spill in H2:H10 and must not be left stale).
# NOTE: single-quoted strings -- "$H$10" would be PowerShell variable
# interpolation ($H, then $1, then 0) and silently mangle the formula.
$wsE.Range("K4").Formula = '=$H$10=$J$10'
$wsE.Range("K7").Formula = '=$H$5>=$J$5'
$wsE.Range("K10").Formula = '=$H$8<=$J$8'

$excel.Calculate()

# ---------------------------------------------------------------------
# Final on-screen state: "2.3E - 8" ends up the active/visible tab, with
# the selections below left on each sheet.
# ---------------------------------------------------------------------
$wsA.Activate()
$wsA.Range("D19").Select()

$ws9B.Activate()
$ws9B.Range("A1:T8").Select()

$wsE.Activate()
$wsE.Range("H8").Select()
